# Task List update: rework the task breakdown (rows 10-20 of TaskList sheet)
# - "Interactive objects" (row 12) is replaced by "Chest and key"
# - "n/a" Act. placeholders are replaced with real actual-hours numbers
# - Est./Act. hours are updated for several existing tasks
# - Four brand-new tasks are appended (Movable objects, Victory screen,
#   Splashscreen, Doors) with their own No./Category/Est./Act.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TaskList")

# Row 10 - Map(level designs, i, ) : Act. n/a -> 20
$ws.Range("F10").Value = 20

# Row 11 - Transitions between maps : Act. n/a -> 4
$ws.Range("F11").Value = 4

# Row 12 - Interactive objects -> Chest and key, Act. n/a -> 5
$ws.Range("B12").Value = "Chest and key"
$ws.Range("F12").Value = 5

# Row 13 - Activity feed : Est. 1 -> 10, Act. 1 -> 5
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 5

# Row 14 - Control edits : unchanged (Est. 2, Act. 2)

# Row 15 - Inventory : Act. n/a -> 10
$ws.Range("F15").Value = 10

# Row 16 - Collision : Est. 3 -> 4
$ws.Range("E16").Value = 4

# Row 17 - new task: Movable objects
$ws.Range("A17").Value = 8
$ws.Range("B17").Value = "Movable objects"
$ws.Range("D17").Value = "SY"
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 6

# Row 18 - new task: Victory screen
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = "Victory screen"
$ws.Range("D18").Value = "DP"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 1

# Row 19 - new task: Splashscreen
$ws.Range("A19").Value = 10
$ws.Range("B19").Value = "Splashscreen"
$ws.Range("D19").Value = "AK"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 1

# Row 20 - new task: Doors
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = "Doors"
$ws.Range("D20").Value = "RC"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 3

# Leave the selection where the author ended up after the edit
$ws.Range("F22").Select() | Out-Null
